# Update the SmartArt ("Req 3 - Support User Submissions..." diagram) text on
# slide 4: append "along with park reviews" to the reviews bullet, and drop
# "and park descriptions" from the user-images bullet.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(5)          # "Content Placeholder 2" graphicFrame (SmartArt)
$sa = $sh.SmartArt

$node1 = $sa.AllNodes.Item(1)
$node1.TextFrame2.TextRange.Text = "We allow the user to submit reviews of the game and how they felt about it, along with park reviews"

$node2 = $sa.AllNodes.Item(2)
$node2.TextFrame2.TextRange.Text = "We did not include user images because we could not figure out how to create a filter to get around possibly inappropriate posts, so we just decided to not include it in the end"
